$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '36.277.32'
$ws.Range('E2').Value = '  -2.65%  '

# Row 3
$ws.Range('D3').Value = '1.985.33'
$ws.Range('E3').Value = '  -2.41%  '

# Row 4
$ws.Range('E4').Value = '  +0.25%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.22'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.12%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.631'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.14%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '62.72'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.06%  '

# Row 8
$ws.Range('E8').Value = '  +0.23%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.378'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.17%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '56.51'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.18%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0802'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +6.93%  '

# Row 12
$ws.Range('E12').Value = '  -0.57%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.866'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.26%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '22.55'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +11.01%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.09'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -6.53%  '

# Row 16
$ws.Range('D16').Value = '2.277.46'
$ws.Range('E16').Value = '  -2.21%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.46'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.04%  '

# Row 18
$ws.Range('D18').Value = '1.993.16'
$ws.Range('E18').Value = '  -1.89%  '

# Row 19
$ws.Range('D19').Value = '36.198.55'
$ws.Range('E19').Value = '  -2.15%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.31'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.24%  '

# Row 21
$ws.Range('D21').Value = '0.0₃0876'
$ws.Range('E21').Value = '  +0.56%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.28'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.64%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '237.40'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.63%  '

# Row 24
$ws.Range('E24').Value = '  -0.14%  '

# Row 25
$ws.Range('E25').Value = '  -9.75%  '

# Row 26
$ws.Range('E26').Value = '  -0.56%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.80'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.97%  '

# Row 28
$ws.Range('E28').Value = '  +23.06%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '159.63'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.58%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.95'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.65%  '

# Row 31
$ws.Range('E31').Value = '  -0.78%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.92'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.96%  '

# Row 33
$ws.Range('E33').Value = '  -5.57%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0624'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.70%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.40'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -6.05%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.37'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +6.32%  '

# Row 37
$ws.Range('E37').Value = '  -6.31%  '

# Row 38
$ws.Range('E38').Value = '  +0.28%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.14'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +14.50%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0996'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.95%  '

# Row 42
$ws.Range('E42').Value = '  +0.61%  '

# Row 43
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0214'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.64%  '

# Row 44
$ws.Range('B44').Value = 'HuobiToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.85'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.61%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.10'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.21%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '93.28'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.89%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '16.31'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.45%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.55'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -6.44%  '

# Row 49
$ws.Range('D49').Value = '1.353.70'
$ws.Range('E49').Value = '  -5.25%  '

# Row 50
$ws.Range('E50').Value = '  -2.62%  '

# Row 51
$ws.Range('D51').Value = '2.172.10'
$ws.Range('E51').Value = '  -1.85%  '
